$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.971.96"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "2.967.88"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'562.96"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'137.17"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.516"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "2.962.24"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").Value = "'5.28"
$ws.Range("E11").Value = "  +8.98%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").Value = "'33.45"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "3.453.94"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'7.04"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "2.964.70"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "59.004.60"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").Value = "'433.85"
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").Value = "'13.52"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'0.717"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").Value = "'6.98"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "'13.05"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").Value = "'79.66"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +8.11%  "
$ws.Range("D29").Value = "'2.53"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").Value = "'7.66"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'25.57"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").Value = "'6.16"
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").Value = "'0.104"
$ws.Range("E33").Value = "  +6.69%  "
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").Value = "0.0₃0760"
$ws.Range("E35").Value = "  +7.15%  "
$ws.Range("D36").Value = "'0.976"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").Value = "'48.32"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'8.70"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").Value = "'395.83"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "2.719.61"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'122.13"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("D48").Value = "'33.96"
$ws.Range("E48").Value = "  +15.26%  "
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").Value = "'1.96"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "'23.05"
$ws.Range("E51").Value = "  +0.25%  "
